# The commit swaps the two theme parts in the deck:
#   ppt/theme/theme1.xml (bound to the Slide Master) used to be the
#   "Integral" / "Red Violet" colour scheme and becomes the stock
#   "Office Theme" colour scheme (the colours that used to live in
#   ppt/theme/theme2.xml, which is bound to the Notes Master).
#
# PowerPoint's object model edits theme colours through
# ThemeColorScheme.Colors(index).RGB (index 1-12, order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), so recreate the
# "Office Theme" palette on the Slide Master's theme that way.

function Get-VbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeThemeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $colorScheme.Colors($i).RGB = Get-VbaRgb $officeThemeHex[$i - 1]
}
